# Adds two new columns ("LabourRate" / "changeLabourRate") to the France
# testdata sheet, just before the existing StageProbability* columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank columns at Z (26) and AA (27); existing Z/AA/AB
# shift right to AB/AC/AD. Insert() without args shifts cells right and
# copies formatting from the column to the left, which lines up with the
# header (style 1) / data (style 8) styling already used by the sheet.
$ws.Columns.Item(26).Insert()
$ws.Columns.Item(27).Insert()

# Populate column Z first (header then data) and then column AA, so the
# shared-string table receives the four new strings in that order.
$ws.Range("Z1").Value = "LabourRate"
$ws.Range("Z2").Value = "61.9"
$ws.Range("AA1").Value = "changeLabourRate"
$ws.Range("AA2").Value = "61"

# Restore explicit column widths for the two new columns.
$ws.Columns.Item(26).ColumnWidth = 11.17
$ws.Columns.Item(27).ColumnWidth = 13
